# Apply the "Save the excel report" edits to Sheet1:
#  - Update the E7/F7 accuracy values (G7's shared formula recalculates automatically)
#  - Move the active cell selection to L8

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E7").Value = 0.67800000000000005
$ws.Range("F7").Value = 0.67300000000000004

$ws.Range("L8").Select()

$wb.Save()
